# Change swt/LED pinouts to accept interrupts, add ISRs to code
# Adds two new BOM line items (SH JST substitute test housing, ZH JST test
# housing) as rows 14-15, plus a trailing formula-only row 16, to the BOM
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pick up the same cell formatting the rest of the BOM table uses, by
# copying format-only from row 12 (the last existing line item). Column E
# intentionally inherits the C/D (not E) format, matching upstream. --------
$ws.Range("B12:D12").Copy() | Out-Null
$ws.Range("B14:D14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B15:D15").PasteSpecial(-4122) | Out-Null

$ws.Range("C12").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null

$ws.Range("F12").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null

$ws.Range("G12").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null

$ws.Range("H12").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 14: SH JST substitute test housing --------------------------------
$ws.Range("B14").Value = "SH JST substitute test housing"
$ws.Range("C14").Value = "TE "
$ws.Range("D14").Value = "1470364-3"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.43
$ws.Range("H14").Formula = "=PRODUCT(E14*F14)"

# --- Row 15: ZH JST test housing -------------------------------------------
$ws.Range("B15").Value = "ZH JST test housing"
$ws.Range("C15").Value = "JST "
$ws.Range("D15").Value = "ZHR-3"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.1
$ws.Range("H15").Formula = "=PRODUCT(E15*F15)"

# --- Row 16: trailing shared-formula-only row -------------------------------
$ws.Range("H16").Formula = "=PRODUCT(E16*F16)"

# Vendor-link text + hyperlinks added last (matches original authoring order,
# which keeps the shared-string table in the same sequence as upstream).
$ws.Range("G14").Value = "www.digikey.com/en/products/detail/te-connectivity-amp-connectors/1470364-3/2077839"
$ws.Range("G15").Value = "www.digikey.com/en/products/detail/jst-sales-america-inc/ZHR-3/608602"

$ws.Hyperlinks.Add($ws.Range("G14"), "http://www.digikey.com/en/products/detail/te-connectivity-amp-connectors/1470364-3/2077839") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G15"), "http://www.digikey.com/en/products/detail/jst-sales-america-inc/ZHR-3/608602") | Out-Null

# Re-apply the thin box border on the two new vendor-link cells (adding the
# hyperlink resets their style to the bare "Hyperlink" named style).
foreach ($cellRef in "G14", "G15") {
  foreach ($edgeIdx in 7, 8, 9, 10) {
    $ws.Range($cellRef).Borders.Item($edgeIdx).LineStyle = 1
    $ws.Range($cellRef).Borders.Item($edgeIdx).Weight = 2
  }
}

# --- View state: drop the old scroll/selection, select G16 ----------------
$ws.Range("G16").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
